$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B4").Value = "inf"
$ws.Range("B6").Value = 466233.8611992272
$ws.Range("B7").Value = 1781550.115583401
$ws.Range("B8").Value = 24017313.31977735
$ws.Range("B10").Value = 4863996.423059861

# --- Sheet: Costs and Revenues ---
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("C2").Value = 128636.3845232744
$ws.Range("D2").Value = 132174.54662458
$ws.Range("F2").Value = 143326.6614602754
$ws.Range("G2").Value = 143326.6614602754
$ws.Range("L2").Value = 143326.6614602754
$ws.Range("P2").Value = 143326.6614602754
$ws.Range("C3").Value = 5404.577564794174
$ws.Range("D3").Value = 19735.29214074445
$ws.Range("E3").Value = 76639.38691994261
$ws.Range("C4").Value = 25403.49618760574
$ws.Range("D4").Value = 23418.84078372031
$ws.Range("C5").Value = 41242.36126897734
$ws.Range("D5").Value = 41744.72270709982
$ws.Range("C6").Value = 56585.94950189719
$ws.Range("D6").Value = 47275.69099301544
$ws.Range("E6").Value = 42615.09486016769
$ws.Range("F6").Value = 119254.4817801102
$ws.Range("G6").Value = 119254.4817801102
$ws.Range("L6").Value = 119254.4817801102
$ws.Range("P6").Value = 119254.4817801102

# --- Sheet: Installed Capacities ---
$ws = $wb.Worksheets.Item("Installed Capacities")
$ws.Range("C3").Value = 349.3009756411624
$ws.Range("D3").Value = 372.3450783073314

# --- Sheet: Added Capacities ---
$ws = $wb.Worksheets.Item("Added Capacities")
$ws.Range("C3").Value = 5.935169881895904
$ws.Range("D3").Value = 23.04410266616899
$ws.Range("E3").Value = 95.39345374605597

# --- Sheet: PV Dispatch ---
$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("G5").Value = 1.404225027703164
$ws.Range("H5").Value = 14.38101956496503
$ws.Range("I5").Value = 54.13638538052629
$ws.Range("J5").Value = 119.1818439450216
$ws.Range("K5").Value = 178.6226893676965
$ws.Range("L5").Value = 221.5972410592673
$ws.Range("M5").Value = 246.5696278956834
$ws.Range("N5").Value = 250.5593822556451
$ws.Range("O5").Value = 236.5961196364217
$ws.Range("P5").Value = 201.9293142649998
$ws.Range("Q5").Value = 151.6405054603802
$ws.Range("R5").Value = 88.20815039645898
$ws.Range("S5").Value = 31.99877781878589
$ws.Range("T5").Value = 6.146995058770605
$ws.Range("U5").Value = 0.1123380022162531
$ws.Range("G6").Value = 0.7513266268508021
$ws.Range("H6").Value = 7.256233475111696
$ws.Range("I6").Value = 25.8680439507842
$ws.Range("J6").Value = 70.98388977628566
$ws.Range("K6").Value = 121.3227737752894
$ws.Range("L6").Value = 163.1334415010919
$ws.Range("M6").Value = 190.3690317244335
$ws.Range("O6").Value = 178.7597172226994
$ws.Range("P6").Value = 143.4704328062729
$ws.Range("Q6").Value = 95.90618485906029
$ws.Range("R6").Value = 46.64815670921035
$ws.Range("S6").Value = 13.95556256453134
$ws.Range("T6").Value = 3.028373552964416
$ws.Range("U6").Value = 0.04942938334544753
$ws.Range("G7").Value = 0.629887005254555
$ws.Range("H7").Value = 5.60026810126323
$ws.Range("I7").Value = 18.94242012165517
$ws.Range("J7").Value = 44.53301127149704
$ws.Range("K7").Value = 73.18141751957465
$ws.Range("L7").Value = 93.64701894484541
$ws.Range("M7").Value = 98.73765119640265
$ws.Range("N7").Value = 96.38989054045393
$ws.Range("O7").Value = 89.03166506998024
$ws.Range("P7").Value = 76.18197016278725
$ws.Range("Q7").Value = 52.74444732181552
$ws.Range("R7").Value = 28.32201025444571
$ws.Range("S7").Value = 10.9772126279362
$ws.Range("T7").Value = 2.691335386087644
$ws.Range("U7").Value = 0.03435747301388486
$ws.Range("G8").Value = 1.496864636411381
$ws.Range("H8").Value = 15.32976495764806
$ws.Range("I8").Value = 57.70787389524983
$ws.Range("J8").Value = 127.0445149346206
$ws.Range("K8").Value = 190.4067949939144
$ws.Range("L8").Value = 236.2164661104893
$ws.Range("M8").Value = 262.8363325882701
$ws.Range("N8").Value = 267.089299236474
$ws.Range("O8").Value = 252.2048515081583
$ws.Range("P8").Value = 215.2510057967523
$ws.Range("Q8").Value = 161.6445410052697
$ws.Range("R8").Value = 94.02742321697652
$ws.Range("S8").Value = 34.10980290222439
$ws.Range("T8").Value = 6.552524945890824
$ws.Range("U8").Value = 0.1197491709129105
$ws.Range("G9").Value = 0.8008931873025619
$ws.Range("H9").Value = 7.734942098422112
$ws.Range("I9").Value = 27.57461193125049
$ws.Range("J9").Value = 75.66684275300214
$ws.Range("K9").Value = 129.3266863022738
$ws.Range("L9").Value = 173.8956896377778
$ws.Range("M9").Value = 202.9280676774956
$ws.Range("O9").Value = 190.5528628579038
$ws.Range("P9").Value = 152.9354718806094
$ws.Range("Q9").Value = 102.2333128212884
$ws.Range("R9").Value = 49.72563140111873
$ws.Range("S9").Value = 14.87623968520328
$ws.Range("T9").Value = 3.228161575136202
$ws.Range("U9").Value = 0.05269034126990541
$ws.Range("G10").Value = 0.6714419444886303
$ws.Range("H10").Value = 5.969729288271644
$ws.Range("I10").Value = 20.19209047607627
$ws.Range("J10").Value = 47.47094547534616
$ws.Range("K10").Value = 78.00934591422448
$ws.Range("L10").Value = 99.82510509242783
$ws.Range("M10").Value = 105.2515768074314
$ws.Range("N10").Value = 102.748929559792
$ws.Range("O10").Value = 94.90526684462934
$ws.Range("P10").Value = 81.20785117706122
$ws.Range("Q10").Value = 56.22410682440704
$ws.Range("R10").Value = 30.19047143127968
$ws.Range("S10").Value = 11.70140188713367
$ws.Range("T10").Value = 2.868888308269602
$ws.Range("U10").Value = 0.03662410606301624
$ws.Range("Q11").Value = 203.0572840332874
$ws.Range("T11").Value = 8.231257986185739
$ws.Range("U12").Value = 0.06618941491321523
$ws.Range("H13").Value = 7.499152202429723
$ws.Range("K13").Value = 97.9950563875785

# --- Sheet: Fed-in Capacity ---
$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("M5").Value = 90.31173819920292
$ws.Range("N5").Value = 84.95722325436918
$ws.Range("O5").Value = 99.00804712831379
$ws.Range("L6").Value = 39.50233892697406
$ws.Range("M6").Value = 17.50199288651834
$ws.Range("O6").Value = 29.7872902773006
$ws.Range("P6").Value = 52.46713803718507
$ws.Range("L8").Value = 108.5919157831171
$ws.Range("M8").Value = 74.04503350661622
$ws.Range("N8").Value = 68.42730627354024
$ws.Range("O8").Value = 84.31378269605858
$ws.Range("L9").Value = 28.74009079028824
$ws.Range("M9").Value = 4.94295693345623
$ws.Range("O9").Value = 17.99414464209616
$ws.Range("P9").Value = 43.00209896284861
$ws.Range("K11").Value = 82.693084352536

# --- Sheet: Unmet Demand ---
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("G5").Value = 422.1859965465589
$ws.Range("H5").Value = 349.7154090116571
$ws.Range("I5").Value = 253.493765682094
$ws.Range("J5").Value = 101.632671480518
$ws.Range("K5").Value = 58.98091051708468
$ws.Range("L5").Value = 14.61922505122197
$ws.Range("O5").Value = 0.9144674394813705
$ws.Range("P5").Value = 52.5812779024964
$ws.Range("Q5").Value = 118.279668409427
$ws.Range("R5").Value = 209.9994294897338
$ws.Range("S5").Value = 235.3909722275223
$ws.Range("T5").Value = 219.367403308014
$ws.Range("U5").Value = 248.7892675345237
$ws.Range("G6").Value = 161.0178692624737
$ws.Range("H6").Value = 137.7783919771127
$ws.Range("I6").Value = 117.0507425405837
$ws.Range("J6").Value = 81.7324700352081
$ws.Range("K6").Value = 11.50936910144016
$ws.Range("Q6").Value = 54.19706493042543
$ws.Range("R6").Value = 154.5724836871435
$ws.Range("S6").Value = 208.0764167058849
$ws.Range("T6").Value = 230.1069956930864
$ws.Range("U6").Value = 249.6497454524407
$ws.Range("G7").Value = 169.2368152776922
$ws.Range("H7").Value = 167.5992636288474
$ws.Range("I7").Value = 167.8087252079146
$ws.Range("J7").Value = 132.4523698397117
$ws.Range("K7").Value = 87.74709330522369
$ws.Range("L7").Value = 61.37839771566142
$ws.Range("M7").Value = 58.02879599051386
$ws.Range("N7").Value = 45.67761186038953
$ws.Range("O7").Value = 71.07018916071037
$ws.Range("P7").Value = 90.14352929744044
$ws.Range("Q7").Value = 151.8424483830337
$ws.Range("R7").Value = 221.1628897093868
$ws.Range("S7").Value = 243.7344203067059
$ws.Range("T7").Value = 218.0970376225092
$ws.Range("U7").Value = 291.2197127438696
$ws.Range("G8").Value = 422.0933569378507
$ws.Range("H8").Value = 348.766663618974
$ws.Range("I8").Value = 249.9222771673705
$ws.Range("J8").Value = 93.77000049091893
$ws.Range("K8").Value = 47.19680489086682
$ws.Range("P8").Value = 39.2595863707439
$ws.Range("Q8").Value = 108.2756328645375
$ws.Range("R8").Value = 204.1801566692163
$ws.Range("S8").Value = 233.2799471440838
$ws.Range("T8").Value = 218.9618734208938
$ws.Range("U8").Value = 248.7818563658271
$ws.Range("G9").Value = 160.968302702022
$ws.Range("H9").Value = 137.2996833538022
$ws.Range("I9").Value = 115.3441745601174
$ws.Range("J9").Value = 77.04951705849162
$ws.Range("K9").Value = 3.505456574455792
$ws.Range("Q9").Value = 47.86993696819729
$ws.Range("R9").Value = 151.4950089952352
$ws.Range("S9").Value = 207.155739585213
$ws.Range("T9").Value = 229.9072076709146
$ws.Range("U9").Value = 249.6464844945162
$ws.Range("G10").Value = 169.1952603384582
$ws.Range("H10").Value = 167.229802441839
$ws.Range("I10").Value = 166.5590548534935
$ws.Range("J10").Value = 129.5144356358626
$ws.Range("K10").Value = 82.91916491057385
$ws.Range("L10").Value = 55.200311568079
$ws.Range("M10").Value = 51.51487037948515
$ws.Range("N10").Value = 39.31857284105146
$ws.Range("O10").Value = 65.19658738606127
$ws.Range("P10").Value = 85.11764828316647
$ws.Range("Q10").Value = 148.3627888804421
$ws.Range("R10").Value = 219.2944285325528
$ws.Range("S10").Value = 243.0102310475085
$ws.Range("T10").Value = 217.9194847003273
$ws.Range("U10").Value = 291.2174461108205
$ws.Range("Q11").Value = 66.8628898365198
$ws.Range("K13").Value = 62.93345443721984
$ws.Range("L13").Value = 29.62548300147245

# --- Sheet: Household Surplus ---
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B3").Value = 397019.7181437318
$ws.Range("B4").Value = 405256.4199213962
